$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("157").Insert()

$ws.Range("A157").Value = 8
$ws.Range("B157").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C157").Value = 'Coquimbo'
$ws.Range("D157").Value = 44778
$ws.Range("E157").Value = 4
$ws.Range("F157").Value = 100112012
$ws.Range("G157").Value = 'Espinaca'
$ws.Range("H157").Value = 'Sin especificar'
$ws.Range("I157").Value = 'Primera'
$ws.Range("J157").Value = 3000
$ws.Range("K157").Value = 500
$ws.Range("L157").Value = 600
$ws.Range("M157").Value = 550
$ws.Range("N157").Value = '$/atado 300 a 500 gramos'
$ws.Range("O157").Value = 'Provincia del Elquí'
$ws.Range("P157").Value = 1100
$ws.Range("Q157").Value = 0.5
$ws.Range("R157").Value = 'Hortaliza'
